$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1686.6
$ws.Range("I19").Value = 1457.2727
$ws.Range("J19").Value = 1966.8889
$ws.Range("K19").Value = 1457.2727
$ws.Range("L19").Value = 1966.8889
$ws.Range("M19").Value = -1282.2727
$ws.Range("N19").Value = -2316.8889

$ws.Range("H74").Value = 3514.2856
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3514.2856
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3514.2856
$ws.Range("N74").Value = -5386.2856

$ws.Range("H77").Value = 3514.2856
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3514.2856
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 17571.428
$ws.Range("N77").Value = -26931.428

$ws.Range("H111").Value = 2694.8572
$ws.Range("I111").Value = 2630.75
$ws.Range("J111").Value = 2780.3333
$ws.Range("K111").Value = 7892.25
$ws.Range("L111").Value = 8340.999899999999
$ws.Range("M111").Value = -4825.25
$ws.Range("N111").Value = -14474.9999

$ws.Range("H116").Value = 6491.5
$ws.Range("I116").Value = 2456.7144
$ws.Range("J116").Value = 12140.2
$ws.Range("K116").Value = 2456.7144
$ws.Range("L116").Value = 12140.2
$ws.Range("M116").Value = 985.2856000000002
$ws.Range("N116").Value = -19024.2

$ws.Range("H132").Value = 32680944
$ws.Range("I132").Value = 3368244.5
$ws.Range("J132").Value = 1000000000
$ws.Range("K132").Value = 10104733.5
$ws.Range("L132").Value = 3000000000
$ws.Range("M132").Value = -10102203.5
$ws.Range("N132").Value = -3000005060

$ws.Range("H138").Value = 2300795.5
$ws.Range("I138").Value = 894.7538500000001
$ws.Range("J138").Value = 9095957
$ws.Range("K138").Value = 2684.26155
$ws.Range("L138").Value = 27287871
$ws.Range("M138").Value = 2455.73845
$ws.Range("N138").Value = -27298151


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2053.682
$ws.Range("I2").Value = 2045.8572
$ws.Range("J2").Value = 2067.375
$ws.Range("K2").Value = 2045.8572
$ws.Range("L2").Value = 2067.375
$ws.Range("M2").Value = -1932.8572
$ws.Range("N2").Value = -2293.375

$ws.Range("H3").Value = 3741.3333
$ws.Range("I3").Value = 1683.3334
$ws.Range("J3").Value = 5799.3335
$ws.Range("K3").Value = 1683.3334
$ws.Range("L3").Value = 5799.3335
$ws.Range("M3").Value = -1568.3334
$ws.Range("N3").Value = -6029.3335

$ws.Range("H45").Value = 1527.3448
$ws.Range("I45").Value = 1187.0588
$ws.Range("J45").Value = 2009.4166
$ws.Range("K45").Value = 1187.0588
$ws.Range("L45").Value = 2009.4166
$ws.Range("M45").Value = -810.0588
$ws.Range("N45").Value = -2763.4166

$ws.Range("H74").Value = 27130.62
$ws.Range("I74").Value = 35313
$ws.Range("J74").Value = 8877.615
$ws.Range("K74").Value = 35313
$ws.Range("L74").Value = 8877.615
$ws.Range("M74").Value = -34439
$ws.Range("N74").Value = -10625.615

$ws.Range("H77").Value = 27130.62
$ws.Range("I77").Value = 35313
$ws.Range("J77").Value = 8877.615
$ws.Range("K77").Value = 176565
$ws.Range("L77").Value = 44388.075
$ws.Range("M77").Value = -172197
$ws.Range("N77").Value = -53124.075

$ws.Range("H116").Value = 2053.682
$ws.Range("I116").Value = 2045.8572
$ws.Range("J116").Value = 2067.375
$ws.Range("K116").Value = 2045.8572
$ws.Range("L116").Value = 2067.375
$ws.Range("M116").Value = 248.1428000000001
$ws.Range("N116").Value = -6655.375

$ws.Range("H122").Value = 2104
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 2265.6
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 6796.799999999999
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -11696.8


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2053.682
$ws.Range("I3").Value = 2045.8572
$ws.Range("J3").Value = 2067.375
$ws.Range("K3").Value = 2045.8572
$ws.Range("L3").Value = 2067.375
$ws.Range("M3").Value = -1931.8572
$ws.Range("N3").Value = -2295.375

$ws.Range("H8").Value = 1105.8
$ws.Range("I8").Value = 868
$ws.Range("J8").Value = 1462.5
$ws.Range("K8").Value = 868
$ws.Range("L8").Value = 1462.5
$ws.Range("M8").Value = -728
$ws.Range("N8").Value = -1742.5

$ws.Range("H105").Value = 2406.1064
$ws.Range("I105").Value = 2399.3809
$ws.Range("J105").Value = 2411.5386
$ws.Range("K105").Value = 2399.3809
$ws.Range("L105").Value = 2411.5386
$ws.Range("M105").Value = -652.3809000000001
$ws.Range("N105").Value = -5905.5386


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17858766
$ws.Range("I31").Value = 31251124
$ws.Range("J31").Value = 2287.6667
$ws.Range("K31").Value = 31251124
$ws.Range("L31").Value = 2287.6667
$ws.Range("M31").Value = -31250829
$ws.Range("N31").Value = -2877.6667

$ws.Range("H34").Value = 17858766
$ws.Range("I34").Value = 31251124
$ws.Range("J34").Value = 2287.6667
$ws.Range("K34").Value = 31251124
$ws.Range("L34").Value = 2287.6667
$ws.Range("M34").Value = -31250922
$ws.Range("N34").Value = -2691.6667

$ws.Range("H132").Value = 515599.72
$ws.Range("I132").Value = 1157.6066
$ws.Range("J132").Value = 3368415
$ws.Range("K132").Value = 3472.8198
$ws.Range("L132").Value = 10105245
$ws.Range("M132").Value = -942.8198000000002
$ws.Range("N132").Value = -10110305


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1110.3529
$ws.Range("I4").Value = 194
$ws.Range("J4").Value = 1392.3077
$ws.Range("K4").Value = 582
$ws.Range("L4").Value = 4176.9231
$ws.Range("M4").Value = -470
$ws.Range("N4").Value = -4400.9231

$ws.Range("H7").Value = 18181996
$ws.Range("I7").Value = 186
$ws.Range("J7").Value = 25000176
$ws.Range("K7").Value = 558
$ws.Range("L7").Value = 75000528
$ws.Range("M7").Value = -446
$ws.Range("N7").Value = -75000752

$ws.Range("H14").Value = 64.666664
$ws.Range("I14").Value = 64.666664
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 193.999992
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -20.99999199999999

$ws.Range("H64").Value = 4743.273
$ws.Range("I64").Value = 3696.4
$ws.Range("J64").Value = 5615.6665
$ws.Range("K64").Value = 11089.2
$ws.Range("L64").Value = 16846.9995
$ws.Range("M64").Value = -10819.2
$ws.Range("N64").Value = -17386.9995

$ws.Range("H67").Value = 4743.273
$ws.Range("I67").Value = 3696.4
$ws.Range("J67").Value = 5615.6665
$ws.Range("K67").Value = 11089.2
$ws.Range("L67").Value = 16846.9995
$ws.Range("M67").Value = -10153.2
$ws.Range("N67").Value = -18718.9995

$ws.Range("H75").Value = 3917.9333
$ws.Range("I75").Value = 2000
$ws.Range("J75").Value = 4054.9285
$ws.Range("K75").Value = 6000
$ws.Range("L75").Value = 12164.7855
$ws.Range("M75").Value = -5002
$ws.Range("N75").Value = -14160.7855

$ws.Range("H78").Value = 3917.9333
$ws.Range("I78").Value = 2000
$ws.Range("J78").Value = 4054.9285
$ws.Range("K78").Value = 18000
$ws.Range("L78").Value = 36494.3565
$ws.Range("M78").Value = -13008
$ws.Range("N78").Value = -46478.3565

$ws.Range("H105").Value = 219350
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 219350
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 658050
$ws.Range("N105").Value = -663292

$ws.Range("H114").Value = 17356166
$ws.Range("I114").Value = 16666956
$ws.Range("J114").Value = 18183218
$ws.Range("K114").Value = 50000868
$ws.Range("L114").Value = 54549654
$ws.Range("M114").Value = -49997614
$ws.Range("N114").Value = -54556162

$ws.Range("H117").Value = 937
$ws.Range("I117").Value = 645
$ws.Range("J117").Value = 1083
$ws.Range("K117").Value = 1935
$ws.Range("L117").Value = 3249
$ws.Range("M117").Value = 1507
$ws.Range("N117").Value = -10133

$ws.Range("H129").Value = 4044.4614
$ws.Range("I129").Value = 4009
$ws.Range("J129").Value = 4066.625
$ws.Range("K129").Value = 12027
$ws.Range("L129").Value = 12199.875
$ws.Range("M129").Value = -7027
$ws.Range("N129").Value = -22199.875

$ws.Range("H137").Value = 16116867
$ws.Range("I137").Value = 2456.3635
$ws.Range("J137").Value = 21834884
$ws.Range("K137").Value = 7369.0905
$ws.Range("L137").Value = 65504652
$ws.Range("M137").Value = -2269.0905
$ws.Range("N137").Value = -65514852

$ws.Range("H140").Value = 3276.3
$ws.Range("I140").Value = 704.2857
$ws.Range("J140").Value = 9277.666999999999
$ws.Range("K140").Value = 2112.8571
$ws.Range("L140").Value = 27833.001
$ws.Range("M140").Value = 3067.1429
$ws.Range("N140").Value = -38193.001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1833.3334
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6840


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2822.889
$ws.Range("I40").Value = 2677
$ws.Range("J40").Value = 2939.6
$ws.Range("K40").Value = 2677
$ws.Range("L40").Value = 2939.6
$ws.Range("M40").Value = -2541
$ws.Range("N40").Value = -3211.6

$ws.Range("H61").Value = 1682.1818
$ws.Range("I61").Value = 1300.6666
$ws.Range("J61").Value = 2140
$ws.Range("K61").Value = 1300.6666
$ws.Range("L61").Value = 2140
$ws.Range("M61").Value = -1098.6666
$ws.Range("N61").Value = -2544

$ws.Range("H113").Value = 1682.1818
$ws.Range("I113").Value = 1300.6666
$ws.Range("J113").Value = 2140
$ws.Range("K113").Value = 1300.6666
$ws.Range("L113").Value = 2140
$ws.Range("M113").Value = 869.3334
$ws.Range("N113").Value = -6480

$ws.Range("H122").Value = 5393.7144
$ws.Range("I122").Value = 5851.2
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 17553.6
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -15103.6
$ws.Range("N122").Value = -17650

$ws.Range("H137").Value = 72130
$ws.Range("I137").Value = 22695
$ws.Range("J137").Value = 171000
$ws.Range("K137").Value = 22695
$ws.Range("L137").Value = 171000
$ws.Range("M137").Value = -17595
$ws.Range("N137").Value = -181200

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2180.1072
$ws.Range("I132").Value = 2113.6978
$ws.Range("J132").Value = 2399.7693
$ws.Range("K132").Value = 6341.0934
$ws.Range("L132").Value = 7199.3079
$ws.Range("M132").Value = -3811.0934
$ws.Range("N132").Value = -12259.3079

$ws.Range("H136").Value = 1162.4138
$ws.Range("I136").Value = 771.3555
$ws.Range("J136").Value = 2516.077
$ws.Range("K136").Value = 2314.0665
$ws.Range("L136").Value = 7548.231000000001
$ws.Range("M136").Value = 235.9335000000001
$ws.Range("N136").Value = -12648.231

